$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value = 10.00636840280725
$ws.Range("B3").Value = 0.18
$ws.Range("B4").Value = 1851.692682370092
$ws.Range("B5").Value = 22457.07626537552
$ws.Range("B8").Value = 0.9385700862638099
$ws.Range("B9").Value = 0.7179030094017506
$ws.Range("B10").Value = 2.814137998255945
$ws.Range("B11").Value = 0.3540906443703401
